$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1619789639338207, 0.3076492955603238, 0.1317988477936555, 0.3630411103355314, 0.3381692893945722, 13)
    3  = @(0.1981620169181928, 0.3116209034719995, 0.1296824127478638, 0.3601144439589501, 0.3140600309705328, 12)
    4  = @(0.2201334750512514, 0.2700445318363022, 0.09544857498369073, 0.3089475278808535, 0.2273517340160526, 11)
    5  = @(0.2251408850412157, 0.2596564006926809, 0.09977641923863516, 0.3158740559758512, 0.2335427658231482, 10)
    6  = @(0.1781092598615317, 0.1781092598615317, 0.05114199694748629, 0.2261459638098507, 0.147805529536099, 9)
    7  = @(0.2144878488911046, 0.2313407205340803, 0.09456382460910477, 0.3075123161909207, 0.235575374136654, 8)
    8  = @(0.3427370110204659, 0.3605818570417189, 0.2023703037307001, 0.4498558699524772, 0.3114971130477678, 8)
    9  = @(0.2987280035122604, 0.2987280035122604, 0.1585806734357395, 0.3982218896993729, 0.2844280147574629, 7)
    10 = @(0.07168159692063568, 0.1326419446783583, 0.02840911087054501, 0.1685500248310424, 0.1705537288270639, 5)
    11 = @(0.018797990165231, 0.2195623404017969, 0.0769818559669036, 0.2774560433057885, 0.3196424909856875, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
